$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 38.33049633333334
$ws.Cells.Item(2, 8).Value = 114.991489
$ws.Cells.Item(2, 9).Value = 0.5317874798120843
$ws.Cells.Item(2, 10).Value = 0.5317874798120843
$ws.Cells.Item(2, 13).Value = 1.070709333333333
$ws.Cells.Item(2, 14).Value = 3.212128
$ws.Cells.Item(2, 15).Value = 0.1056496288760636
$ws.Cells.Item(2, 16).Value = 0.1056496288760636
$ws.Cells.Item(2, 17).Value = 41.04082017539911
$ws.Cells.Item(2, 18).Value = 369.367381578592
$ws.Cells.Item(2, 19).Value = 0.05618314988308389
$ws.Cells.Item(2, 20).Value = 0.05618314988308389
$ws.Cells.Item(3, 7).Value = 38.33049633333334
$ws.Cells.Item(3, 8).Value = 114.991489
$ws.Cells.Item(3, 9).Value = 0.5317874798120843
$ws.Cells.Item(3, 10).Value = 0.5317874798120843
$ws.Cells.Item(3, 15).Value = 0.152405028485123
$ws.Cells.Item(3, 16).Value = 0.152405028485123
$ws.Cells.Item(3, 17).Value = 59.20349587997111
$ws.Cells.Item(3, 18).Value = 532.83146291974
$ws.Cells.Item(3, 19).Value = 0.08104708600879246
$ws.Cells.Item(3, 20).Value = 0.08104708600879246
$ws.Cells.Item(4, 7).Value = 38.33049633333334
$ws.Cells.Item(4, 8).Value = 114.991489
$ws.Cells.Item(4, 9).Value = 0.5317874798120843
$ws.Cells.Item(4, 10).Value = 0.5317874798120843
$ws.Cells.Item(4, 13).Value = 3.115442666666667
$ws.Cells.Item(4, 14).Value = 9.346328
$ws.Cells.Item(4, 15).Value = 0.3074086974597408
$ws.Cells.Item(4, 16).Value = 0.3074086974597408
$ws.Cells.Item(4, 17).Value = 119.4164637113769
$ws.Cells.Item(4, 18).Value = 1074.748173402392
$ws.Cells.Item(4, 19).Value = 0.163476096494431
$ws.Cells.Item(4, 20).Value = 0.163476096494431
$ws.Cells.Item(5, 7).Value = 38.33049633333334
$ws.Cells.Item(5, 8).Value = 114.991489
$ws.Cells.Item(5, 9).Value = 0.5317874798120843
$ws.Cells.Item(5, 10).Value = 0.5317874798120843
$ws.Cells.Item(5, 13).Value = 1.972894
$ws.Cells.Item(5, 14).Value = 5.918682
$ws.Cells.Item(5, 15).Value = 0.1946704977931882
$ws.Cells.Item(5, 16).Value = 0.1946704977931882
$ws.Cells.Item(5, 17).Value = 75.62200623305533
$ws.Cells.Item(5, 18).Value = 680.598056097498
$ws.Cells.Item(5, 19).Value = 0.1035233334152035
$ws.Cells.Item(5, 20).Value = 0.1035233334152035
$ws.Cells.Item(6, 7).Value = 38.33049633333334
$ws.Cells.Item(6, 8).Value = 114.991489
$ws.Cells.Item(6, 9).Value = 0.5317874798120843
$ws.Cells.Item(6, 10).Value = 0.5317874798120843
$ws.Cells.Item(6, 13).Value = 2.430930666666667
$ws.Cells.Item(6, 14).Value = 7.292792
$ws.Cells.Item(6, 15).Value = 0.2398661473858844
$ws.Cells.Item(6, 16).Value = 0.2398661473858844
$ws.Cells.Item(6, 17).Value = 93.17877900525423
$ws.Cells.Item(6, 18).Value = 838.609011047288
$ws.Cells.Item(6, 19).Value = 0.1275578140105734
$ws.Cells.Item(6, 20).Value = 0.1275578140105734
$ws.Cells.Item(7, 9).Value = 0.2073226210890634
$ws.Cells.Item(7, 10).Value = 0.2073226210890634
$ws.Cells.Item(7, 13).Value = 1.070709333333333
$ws.Cells.Item(7, 14).Value = 3.212128
$ws.Cells.Item(7, 15).Value = 0.1056496288760636
$ws.Cells.Item(7, 16).Value = 0.1056496288760636
$ws.Cells.Item(7, 17).Value = 16.00017061969067
$ws.Cells.Item(7, 18).Value = 144.001535577216
$ws.Cells.Item(7, 19).Value = 0.02190355797567231
$ws.Cells.Item(7, 20).Value = 0.02190355797567231
$ws.Cells.Item(8, 9).Value = 0.2073226210890634
$ws.Cells.Item(8, 10).Value = 0.2073226210890634
$ws.Cells.Item(8, 15).Value = 0.152405028485123
$ws.Cells.Item(8, 16).Value = 0.152405028485123
$ws.Cells.Item(8, 19).Value = 0.03159700997268906
$ws.Cells.Item(8, 20).Value = 0.03159700997268906
$ws.Cells.Item(9, 9).Value = 0.2073226210890634
$ws.Cells.Item(9, 10).Value = 0.2073226210890634
$ws.Cells.Item(9, 13).Value = 3.115442666666667
$ws.Cells.Item(9, 14).Value = 9.346328
$ws.Cells.Item(9, 15).Value = 0.3074086974597408
$ws.Cells.Item(9, 16).Value = 0.3074086974597408
$ws.Cells.Item(9, 17).Value = 46.55569225995734
$ws.Cells.Item(9, 18).Value = 419.001230339616
$ws.Cells.Item(9, 19).Value = 0.06373277690292835
$ws.Cells.Item(9, 20).Value = 0.06373277690292835
$ws.Cells.Item(10, 9).Value = 0.2073226210890634
$ws.Cells.Item(10, 10).Value = 0.2073226210890634
$ws.Cells.Item(10, 13).Value = 1.972894
$ws.Cells.Item(10, 14).Value = 5.918682
$ws.Cells.Item(10, 15).Value = 0.1946704977931882
$ws.Cells.Item(10, 16).Value = 0.1946704977931882
$ws.Cells.Item(10, 17).Value = 29.481988838456
$ws.Cells.Item(10, 18).Value = 265.337899546104
$ws.Cells.Item(10, 19).Value = 0.04035959785119651
$ws.Cells.Item(10, 20).Value = 0.0403595978511965
$ws.Cells.Item(11, 9).Value = 0.2073226210890634
$ws.Cells.Item(11, 10).Value = 0.2073226210890634
$ws.Cells.Item(11, 13).Value = 2.430930666666667
$ws.Cells.Item(11, 14).Value = 7.292792
$ws.Cells.Item(11, 15).Value = 0.2398661473858844
$ws.Cells.Item(11, 16).Value = 0.2398661473858844
$ws.Cells.Item(11, 17).Value = 36.32667075966934
$ws.Cells.Item(11, 18).Value = 326.940036837024
$ws.Cells.Item(11, 19).Value = 0.04972967838657713
$ws.Cells.Item(11, 20).Value = 0.04972967838657712
$ws.Cells.Item(12, 7).Value = 8.167063666666666
$ws.Cells.Item(12, 8).Value = 24.501191
$ws.Cells.Item(12, 9).Value = 0.1133077476219524
$ws.Cells.Item(12, 10).Value = 0.1133077476219524
$ws.Cells.Item(12, 13).Value = 1.070709333333333
$ws.Cells.Item(12, 14).Value = 3.212128
$ws.Cells.Item(12, 15).Value = 0.1056496288760636
$ws.Cells.Item(12, 16).Value = 0.1056496288760636
$ws.Cells.Item(12, 17).Value = 8.744551293827554
$ws.Cells.Item(12, 18).Value = 78.70096164444799
$ws.Cells.Item(12, 19).Value = 0.01197092148504196
$ws.Cells.Item(12, 20).Value = 0.01197092148504196
$ws.Cells.Item(13, 7).Value = 8.167063666666666
$ws.Cells.Item(13, 8).Value = 24.501191
$ws.Cells.Item(13, 9).Value = 0.1133077476219524
$ws.Cells.Item(13, 10).Value = 0.1133077476219524
$ws.Cells.Item(13, 15).Value = 0.152405028485123
$ws.Cells.Item(13, 16).Value = 0.152405028485123
$ws.Cells.Item(13, 17).Value = 12.61446540989555
$ws.Cells.Item(13, 18).Value = 113.53018868906
$ws.Cells.Item(13, 19).Value = 0.01726867050390879
$ws.Cells.Item(13, 20).Value = 0.01726867050390879
$ws.Cells.Item(14, 7).Value = 8.167063666666666
$ws.Cells.Item(14, 8).Value = 24.501191
$ws.Cells.Item(14, 9).Value = 0.1133077476219524
$ws.Cells.Item(14, 10).Value = 0.1133077476219524
$ws.Cells.Item(14, 13).Value = 3.115442666666667
$ws.Cells.Item(14, 14).Value = 9.346328
$ws.Cells.Item(14, 15).Value = 0.3074086974597408
$ws.Cells.Item(14, 16).Value = 0.3074086974597408
$ws.Cells.Item(14, 17).Value = 25.44401860851644
$ws.Cells.Item(14, 18).Value = 228.996167476648
$ws.Cells.Item(14, 19).Value = 0.03483178710856143
$ws.Cells.Item(14, 20).Value = 0.03483178710856144
$ws.Cells.Item(15, 7).Value = 8.167063666666666
$ws.Cells.Item(15, 8).Value = 24.501191
$ws.Cells.Item(15, 9).Value = 0.1133077476219524
$ws.Cells.Item(15, 10).Value = 0.1133077476219524
$ws.Cells.Item(15, 13).Value = 1.972894
$ws.Cells.Item(15, 14).Value = 5.918682
$ws.Cells.Item(15, 15).Value = 0.1946704977931882
$ws.Cells.Item(15, 16).Value = 0.1946704977931882
$ws.Cells.Item(15, 17).Value = 16.11275090558467
$ws.Cells.Item(15, 18).Value = 145.014758150262
$ws.Cells.Item(15, 19).Value = 0.02205767563339042
$ws.Cells.Item(15, 20).Value = 0.02205767563339042
$ws.Cells.Item(16, 7).Value = 8.167063666666666
$ws.Cells.Item(16, 8).Value = 24.501191
$ws.Cells.Item(16, 9).Value = 0.1133077476219524
$ws.Cells.Item(16, 10).Value = 0.1133077476219524
$ws.Cells.Item(16, 13).Value = 2.430930666666667
$ws.Cells.Item(16, 14).Value = 7.292792
$ws.Cells.Item(16, 15).Value = 0.2398661473858844
$ws.Cells.Item(16, 16).Value = 0.2398661473858844
$ws.Cells.Item(16, 17).Value = 19.85356552391911
$ws.Cells.Item(16, 18).Value = 178.682089715272
$ws.Cells.Item(16, 19).Value = 0.02717869289104983
$ws.Cells.Item(16, 20).Value = 0.02717869289104983
$ws.Cells.Item(17, 7).Value = 5.834252333333334
$ws.Cells.Item(17, 8).Value = 17.502757
$ws.Cells.Item(17, 9).Value = 0.08094292121735479
$ws.Cells.Item(17, 10).Value = 0.08094292121735479
$ws.Cells.Item(17, 13).Value = 1.070709333333333
$ws.Cells.Item(17, 14).Value = 3.212128
$ws.Cells.Item(17, 15).Value = 0.1056496288760636
$ws.Cells.Item(17, 16).Value = 0.1056496288760636
$ws.Cells.Item(17, 17).Value = 6.246788426321778
$ws.Cells.Item(17, 18).Value = 56.221095836896
$ws.Cells.Item(17, 19).Value = 0.008551589586757991
$ws.Cells.Item(17, 20).Value = 0.008551589586757989
$ws.Cells.Item(18, 7).Value = 5.834252333333334
$ws.Cells.Item(18, 8).Value = 17.502757
$ws.Cells.Item(18, 9).Value = 0.08094292121735479
$ws.Cells.Item(18, 10).Value = 0.08094292121735479
$ws.Cells.Item(18, 15).Value = 0.152405028485123
$ws.Cells.Item(18, 16).Value = 0.152405028485123
$ws.Cells.Item(18, 17).Value = 9.01131388895778
$ws.Cells.Item(18, 18).Value = 81.10182500062001
$ws.Cells.Item(18, 19).Value = 0.01233610821380002
$ws.Cells.Item(18, 20).Value = 0.01233610821380002
$ws.Cells.Item(19, 7).Value = 5.834252333333334
$ws.Cells.Item(19, 8).Value = 17.502757
$ws.Cells.Item(19, 9).Value = 0.08094292121735479
$ws.Cells.Item(19, 10).Value = 0.08094292121735479
$ws.Cells.Item(19, 13).Value = 3.115442666666667
$ws.Cells.Item(19, 14).Value = 9.346328
$ws.Cells.Item(19, 15).Value = 0.3074086974597408
$ws.Cells.Item(19, 16).Value = 0.3074086974597408
$ws.Cells.Item(19, 17).Value = 18.17627864736622
$ws.Cells.Item(19, 18).Value = 163.586507826296
$ws.Cells.Item(19, 19).Value = 0.02488255798001345
$ws.Cells.Item(19, 20).Value = 0.02488255798001345
$ws.Cells.Item(20, 7).Value = 5.834252333333334
$ws.Cells.Item(20, 8).Value = 17.502757
$ws.Cells.Item(20, 9).Value = 0.08094292121735479
$ws.Cells.Item(20, 10).Value = 0.08094292121735479
$ws.Cells.Item(20, 13).Value = 1.972894
$ws.Cells.Item(20, 14).Value = 5.918682
$ws.Cells.Item(20, 15).Value = 0.1946704977931882
$ws.Cells.Item(20, 16).Value = 0.1946704977931882
$ws.Cells.Item(20, 17).Value = 11.51036142291933
$ws.Cells.Item(20, 18).Value = 103.593252806274
$ws.Cells.Item(20, 19).Value = 0.01575719876621727
$ws.Cells.Item(20, 20).Value = 0.01575719876621727
$ws.Cells.Item(21, 7).Value = 5.834252333333334
$ws.Cells.Item(21, 8).Value = 17.502757
$ws.Cells.Item(21, 9).Value = 0.08094292121735479
$ws.Cells.Item(21, 10).Value = 0.08094292121735479
$ws.Cells.Item(21, 13).Value = 2.430930666666667
$ws.Cells.Item(21, 14).Value = 7.292792
$ws.Cells.Item(21, 15).Value = 0.2398661473858844
$ws.Cells.Item(21, 16).Value = 0.2398661473858844
$ws.Cells.Item(21, 17).Value = 14.18266291417156
$ws.Cells.Item(21, 18).Value = 127.643966227544
$ws.Cells.Item(21, 19).Value = 0.01941546667056605
$ws.Cells.Item(21, 20).Value = 0.01941546667056605
$ws.Cells.Item(22, 7).Value = 4.803262333333334
$ws.Cells.Item(22, 8).Value = 14.409787
$ws.Cells.Item(22, 9).Value = 0.06663923025954499
$ws.Cells.Item(22, 10).Value = 0.066639230259545
$ws.Cells.Item(22, 13).Value = 1.070709333333333
$ws.Cells.Item(22, 14).Value = 3.212128
$ws.Cells.Item(22, 15).Value = 0.1056496288760636
$ws.Cells.Item(22, 16).Value = 0.1056496288760636
$ws.Cells.Item(22, 17).Value = 5.142897810748445
$ws.Cells.Item(22, 18).Value = 46.286080296736
$ws.Cells.Item(22, 19).Value = 0.007040409945507479
$ws.Cells.Item(22, 20).Value = 0.00704040994550748
$ws.Cells.Item(23, 7).Value = 4.803262333333334
$ws.Cells.Item(23, 8).Value = 14.409787
$ws.Cells.Item(23, 9).Value = 0.06663923025954499
$ws.Cells.Item(23, 10).Value = 0.066639230259545
$ws.Cells.Item(23, 15).Value = 0.152405028485123
$ws.Cells.Item(23, 16).Value = 0.152405028485123
$ws.Cells.Item(23, 17).Value = 7.418894847824444
$ws.Cells.Item(23, 18).Value = 66.77005363042001
$ws.Cells.Item(23, 19).Value = 0.01015615378593262
$ws.Cells.Item(23, 20).Value = 0.01015615378593263
$ws.Cells.Item(24, 7).Value = 4.803262333333334
$ws.Cells.Item(24, 8).Value = 14.409787
$ws.Cells.Item(24, 9).Value = 0.06663923025954499
$ws.Cells.Item(24, 10).Value = 0.066639230259545
$ws.Cells.Item(24, 13).Value = 3.115442666666667
$ws.Cells.Item(24, 14).Value = 9.346328
$ws.Cells.Item(24, 15).Value = 0.3074086974597408
$ws.Cells.Item(24, 16).Value = 0.3074086974597408
$ws.Cells.Item(24, 17).Value = 14.96428841245956
$ws.Cells.Item(24, 18).Value = 134.678595712136
$ws.Cells.Item(24, 19).Value = 0.02048547897380647
$ws.Cells.Item(24, 20).Value = 0.02048547897380647
$ws.Cells.Item(25, 7).Value = 4.803262333333334
$ws.Cells.Item(25, 8).Value = 14.409787
$ws.Cells.Item(25, 9).Value = 0.06663923025954499
$ws.Cells.Item(25, 10).Value = 0.066639230259545
$ws.Cells.Item(25, 13).Value = 1.972894
$ws.Cells.Item(25, 14).Value = 5.918682
$ws.Cells.Item(25, 15).Value = 0.1946704977931882
$ws.Cells.Item(25, 16).Value = 0.1946704977931882
$ws.Cells.Item(25, 17).Value = 9.476327437859334
$ws.Cells.Item(25, 18).Value = 85.286946940734
$ws.Cells.Item(25, 19).Value = 0.01297269212718052
$ws.Cells.Item(25, 20).Value = 0.01297269212718052
$ws.Cells.Item(26, 7).Value = 4.803262333333334
$ws.Cells.Item(26, 8).Value = 14.409787
$ws.Cells.Item(26, 9).Value = 0.06663923025954499
$ws.Cells.Item(26, 10).Value = 0.066639230259545
$ws.Cells.Item(26, 13).Value = 2.430930666666667
$ws.Cells.Item(26, 14).Value = 7.292792
$ws.Cells.Item(26, 15).Value = 0.2398661473858844
$ws.Cells.Item(26, 16).Value = 0.2398661473858844
$ws.Cells.Item(26, 17).Value = 11.67639770614489
$ws.Cells.Item(26, 18).Value = 105.087579355304
$ws.Cells.Item(26, 19).Value = 0.0159844954271179
$ws.Cells.Item(26, 20).Value = 0.0159844954271179
